# Update the "想去人数" (interested count) values in column F
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 4146
    4  = 2387
    5  = 476
    9  = 207
    10 = 120
    11 = 107
    12 = 146
    13 = 1551
    14 = 284
    15 = 3081
    16 = 211
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$updates2 = @{
    3  = 4146
    4  = 2387
    5  = 476
    11 = 207
    12 = 120
    13 = 107
    14 = 146
    17 = 1551
    18 = 284
    19 = 3081
    20 = 211
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates2.Keys) {
    $ws4.Range("F$row").Value = $updates2[$row]
}
